# Course date update on the "CustomShape 4" date/time textbox (slide 1, shape 7):
#   "Tue 25th and Wed 26th June 2024"  ->  "Tue 27th and Wed 28th May 2025"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(7)
$tr = $sh.TextFrame.TextRange

# Run 1: "Tue 25" -> "Tue 27" (characters 1-6; replacing the full run keeps it as one run)
$tr.Characters(1, 6).Text = "Tue 27"

# Run 2 ("th", characters 7-8) is unchanged.

# Run 3: " and Wed 26" -> " and Wed " (characters 9-19, shrinks by 2 chars to 9-17)
$tr.Characters(9, 11).Text = " and Wed "

# Insert the new day-of-month run "28" immediately after " and Wed " (now at 9-17).
$wedRun = $tr.Characters(9, 9)
[void]$wedRun.InsertAfter("28")

# Re-fetch the freshly inserted "28" text (now at characters 18-19) and give it its
# own explicit formatting so it stays a distinct run rather than merging back into
# the preceding " and Wed " run.
$newRun = $tr.Characters(18, 2)
$newRun.Font.Name = "Calibri"
$newRun.Font.Size = 24
$newRun.Font.Bold = 0
$newRun.Font.Italic = 0
$newRun.Font.BaselineOffset = 0
$newRun.Font.Color.RGB = 15921906

# The old "th" run (now still at characters 20-21, net offset unchanged since we
# removed 2 chars "26" and inserted 2 chars "28") is left as-is.

# Run 6: " June 2024" -> " May 2025" (characters 22-31)
$tr.Characters(22, 10).Text = " May 2025"
